$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "Storability" good/bad values from column A into column B
$ws.Range("B1:B17").Value2 = $ws.Range("A1:A17").Value2

# Move the header style (applied to the old A1 "Storability" cell) onto the new B1 header cell
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Reset column A's old header cell back to the default (unstyled) format, using A2 as the
# source of the default formatting
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Fill column A with the new sample names (header "Samples" + sample ids), shifting the
# original data into column B
$samples = @(
    "Samples",
    "V1_1_t0",
    "V1_2_t0",
    "V1_3_t0",
    "V1_4_t0",
    "V2_1_t0",
    "V2_2_t0",
    "V2_3_t0",
    "V2_4_t0",
    "V5_1_t0",
    "V5_2_t0",
    "V5_3_t0",
    "V5_4_t0",
    "V6_1_t0",
    "V6_2_t0",
    "V6_3_t0",
    "V6_4_t0"
)

for ($i = 0; $i -lt $samples.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value2 = $samples[$i]
}

$ws.Range("F7").Select()
